$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:2").Delete()
